# Auto-generated edit script: Add data for 2024-10-18
# Applies updated 2024 (column K) values (and a few F/J corrections) across
# the Citywide Totals, By Neighborhood, and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @(
    @{ Sheet = 'Citywide Totals'; Cells = @(@{ Cell = 'K2'; Value = 6034 }, @{ Cell = 'K3'; Value = 6218 }, @{ Cell = 'F4'; Value = 1915 }, @{ Cell = 'J4'; Value = 1836 }, @{ Cell = 'K4'; Value = 1301 }, @{ Cell = 'K5'; Value = 440 }, @{ Cell = 'K6'; Value = 6833 }, @{ Cell = 'F7'; Value = 24108 }, @{ Cell = 'J7'; Value = 29301 }, @{ Cell = 'K7'; Value = 20826 }) },
    @{ Sheet = 'By Neighborhood'; Cells = @(@{ Cell = 'K2'; Value = 183 }, @{ Cell = 'K7'; Value = 608 }, @{ Cell = 'K8'; Value = 1369 }, @{ Cell = 'K9'; Value = 91 }, @{ Cell = 'K11'; Value = 390 }, @{ Cell = 'K18'; Value = 138 }, @{ Cell = 'K19'; Value = 603 }, @{ Cell = 'K20'; Value = 495 }, @{ Cell = 'K23'; Value = 212 }, @{ Cell = 'K24'; Value = 61 }, @{ Cell = 'K25'; Value = 98 }, @{ Cell = 'K27'; Value = 194 }, @{ Cell = 'K33'; Value = 903 }, @{ Cell = 'K34'; Value = 120 }, @{ Cell = 'K37'; Value = 710 }, @{ Cell = 'K38'; Value = 20 }, @{ Cell = 'K39'; Value = 26 }, @{ Cell = 'K41'; Value = 146 }, @{ Cell = 'K42'; Value = 774 }, @{ Cell = 'K43'; Value = 175 }, @{ Cell = 'K47'; Value = 145 }, @{ Cell = 'K48'; Value = 262 }, @{ Cell = 'K51'; Value = 268 }, @{ Cell = 'K52'; Value = 544 }, @{ Cell = 'K53'; Value = 267 }, @{ Cell = 'K55'; Value = 229 }, @{ Cell = 'F63'; Value = 199 }, @{ Cell = 'J63'; Value = 116 }, @{ Cell = 'K63'; Value = 61 }, @{ Cell = 'K65'; Value = 490 }, @{ Cell = 'K67'; Value = 811 }, @{ Cell = 'K72'; Value = 102 }, @{ Cell = 'K73'; Value = 184 }, @{ Cell = 'K76'; Value = 281 }, @{ Cell = 'K77'; Value = 145 }, @{ Cell = 'K82'; Value = 22 }, @{ Cell = 'K83'; Value = 462 }, @{ Cell = 'K85'; Value = 966 }, @{ Cell = 'K88'; Value = 224 }, @{ Cell = 'K89'; Value = 303 }, @{ Cell = 'K90'; Value = 192 }, @{ Cell = 'K91'; Value = 236 }, @{ Cell = 'K93'; Value = 77 }, @{ Cell = 'K94'; Value = 280 }, @{ Cell = 'K95'; Value = 350 }, @{ Cell = 'K96'; Value = 219 }, @{ Cell = 'F101'; Value = 24108 }, @{ Cell = 'J101'; Value = 29301 }, @{ Cell = 'K101'; Value = 20826 }) },
    @{ Sheet = 'West Ridge'; Cells = @(@{ Cell = 'K2'; Value = 68 }, @{ Cell = 'K7'; Value = 219 }) },
    @{ Sheet = 'Auburn Gresham'; Cells = @(@{ Cell = 'K2'; Value = 200 }, @{ Cell = 'K7'; Value = 608 }) },
    @{ Sheet = 'Belmont Cragin'; Cells = @(@{ Cell = 'K2'; Value = 138 }, @{ Cell = 'K3'; Value = 101 }, @{ Cell = 'K6'; Value = 125 }, @{ Cell = 'K7'; Value = 390 }) },
    @{ Sheet = 'Uptown'; Cells = @(@{ Cell = 'K3'; Value = 92 }, @{ Cell = 'K7'; Value = 303 }) },
    @{ Sheet = 'South Shore'; Cells = @(@{ Cell = 'K3'; Value = 335 }, @{ Cell = 'K7'; Value = 966 }) },
    @{ Sheet = 'Little Village'; Cells = @(@{ Cell = 'K2'; Value = 147 }, @{ Cell = 'K4'; Value = 30 }, @{ Cell = 'K7'; Value = 544 }) },
    @{ Sheet = 'Logan Square'; Cells = @(@{ Cell = 'K2'; Value = 69 }, @{ Cell = 'K7'; Value = 267 }) },
    @{ Sheet = 'Austin'; Cells = @(@{ Cell = 'K2'; Value = 375 }, @{ Cell = 'K3'; Value = 418 }, @{ Cell = 'K6'; Value = 463 }, @{ Cell = 'K7'; Value = 1369 }) },
    @{ Sheet = 'South Chicago'; Cells = @(@{ Cell = 'K2'; Value = 159 }, @{ Cell = 'K3'; Value = 165 }, @{ Cell = 'K7'; Value = 462 }) },
    @{ Sheet = 'Garfield Park'; Cells = @(@{ Cell = 'K3'; Value = 330 }, @{ Cell = 'K4'; Value = 45 }, @{ Cell = 'K7'; Value = 903 }) },
    @{ Sheet = 'West Pullman'; Cells = @(@{ Cell = 'K3'; Value = 124 }, @{ Cell = 'K7'; Value = 350 }) },
    @{ Sheet = 'Grand Crossing'; Cells = @(@{ Cell = 'K3'; Value = 235 }, @{ Cell = 'K5'; Value = 31 }, @{ Cell = 'K6'; Value = 207 }, @{ Cell = 'K7'; Value = 710 }) },
    @{ Sheet = 'New City'; Cells = @(@{ Cell = 'K2'; Value = 162 }, @{ Cell = 'K7'; Value = 490 }) },
    @{ Sheet = 'North Lawndale'; Cells = @(@{ Cell = 'K2'; Value = 225 }, @{ Cell = 'K3'; Value = 291 }, @{ Cell = 'K6'; Value = 231 }, @{ Cell = 'K7'; Value = 811 }) },
    @{ Sheet = 'Lake View'; Cells = @(@{ Cell = 'K4'; Value = 37 }, @{ Cell = 'K7'; Value = 262 }) },
    @{ Sheet = 'Chatham'; Cells = @(@{ Cell = 'K4'; Value = 29 }, @{ Cell = 'K6'; Value = 194 }, @{ Cell = 'K7'; Value = 603 }) },
    @{ Sheet = 'River North'; Cells = @(@{ Cell = 'K2'; Value = 62 }, @{ Cell = 'K3'; Value = 54 }, @{ Cell = 'K7'; Value = 281 }) },
    @{ Sheet = 'Hermosa'; Cells = @(@{ Cell = 'K3'; Value = 29 }, @{ Cell = 'K7'; Value = 146 }) },
    @{ Sheet = 'Humboldt Park'; Cells = @(@{ Cell = 'K3'; Value = 233 }, @{ Cell = 'K6'; Value = 287 }, @{ Cell = 'K7'; Value = 774 }) },
    @{ Sheet = 'Lower West Side'; Cells = @(@{ Cell = 'K2'; Value = 72 }, @{ Cell = 'K7'; Value = 229 }) },
    @{ Sheet = 'Dunning'; Cells = @(@{ Cell = 'K2'; Value = 23 }, @{ Cell = 'K7'; Value = 61 }) },
    @{ Sheet = 'Douglas'; Cells = @(@{ Cell = 'K3'; Value = 74 }, @{ Cell = 'K7'; Value = 212 }) },
    @{ Sheet = 'Washington Park'; Cells = @(@{ Cell = 'K2'; Value = 60 }, @{ Cell = 'K3'; Value = 112 }, @{ Cell = 'K7'; Value = 236 }) },
    @{ Sheet = 'Chicago Lawn'; Cells = @(@{ Cell = 'K2'; Value = 164 }, @{ Cell = 'K3'; Value = 161 }, @{ Cell = 'K7'; Value = 495 }) },
    @{ Sheet = 'Calumet Heights'; Cells = @(@{ Cell = 'K3'; Value = 45 }, @{ Cell = 'K7'; Value = 138 }) },
    @{ Sheet = 'West Lawn'; Cells = @(@{ Cell = 'K3'; Value = 18 }, @{ Cell = 'K7'; Value = 77 }) },
    @{ Sheet = 'Garfield Ridge'; Cells = @(@{ Cell = 'K2'; Value = 47 }, @{ Cell = 'K7'; Value = 120 }) },
    @{ Sheet = 'West Loop'; Cells = @(@{ Cell = 'K3'; Value = 56 }, @{ Cell = 'K7'; Value = 280 }) },
    @{ Sheet = 'East Side'; Cells = @(@{ Cell = 'K2'; Value = 35 }, @{ Cell = 'K5'; Value = 3 }, @{ Cell = 'K7'; Value = 98 }) },
    @{ Sheet = 'Kenwood'; Cells = @(@{ Cell = 'K6'; Value = 47 }, @{ Cell = 'K7'; Value = 145 }) },
    @{ Sheet = 'Greektown'; Cells = @(@{ Cell = 'K5'; Value = 16 }, @{ Cell = 'K6'; Value = 26 }) },
    @{ Sheet = 'Avalon Park'; Cells = @(@{ Cell = 'K6'; Value = 23 }, @{ Cell = 'K7'; Value = 91 }) },
    @{ Sheet = 'Portage Park'; Cells = @(@{ Cell = 'K3'; Value = 47 }, @{ Cell = 'K7'; Value = 184 }) },
    @{ Sheet = 'Albany Park'; Cells = @(@{ Cell = 'K3'; Value = 47 }, @{ Cell = 'K7'; Value = 183 }) },
    @{ Sheet = 'United Center'; Cells = @(@{ Cell = 'K2'; Value = 57 }, @{ Cell = 'K6'; Value = 93 }, @{ Cell = 'K7'; Value = 224 }) },
    @{ Sheet = 'Edgewater'; Cells = @(@{ Cell = 'K2'; Value = 53 }, @{ Cell = 'K4'; Value = 24 }, @{ Cell = 'K7'; Value = 194 }) },
    @{ Sheet = 'Washington Heights'; Cells = @(@{ Cell = 'K2'; Value = 71 }, @{ Cell = 'K7'; Value = 192 }) },
    @{ Sheet = 'Little Italy, UIC'; Cells = @(@{ Cell = 'K2'; Value = 76 }, @{ Cell = 'K7'; Value = 268 }) },
    @{ Sheet = 'Hyde Park'; Cells = @(@{ Cell = 'K2'; Value = 35 }, @{ Cell = 'K7'; Value = 175 }) },
    @{ Sheet = 'Old Town'; Cells = @(@{ Cell = 'K6'; Value = 52 }, @{ Cell = 'K7'; Value = 102 }) },
    @{ Sheet = 'Sheffield & DePaul'; Cells = @(@{ Cell = 'K2'; Value = 4 }, @{ Cell = 'K6'; Value = 22 }) },
    @{ Sheet = 'Riverdale'; Cells = @(@{ Cell = 'K4'; Value = 9 }, @{ Cell = 'K7'; Value = 145 }) },
    @{ Sheet = 'Grant Park'; Cells = @(@{ Cell = 'K5'; Value = 10 }, @{ Cell = 'K6'; Value = 20 }) }
)

foreach ($entry in $sheetUpdates) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($c in $entry.Cells) {
        $ws.Range($c.Cell).Value = $c.Value
    }
}

Write-Host "Applied updates to $($sheetUpdates.Count) sheets"
